$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

# --- Fill in release 96 (idx C) and release 97 (idx D) data columns ---
# Match the bold header formatting already used in B1 for the two new columns.
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range("C1").Value = 96
$ws.Range("D1").Value = 97
$ws.Range("C2").Value = 41674
$ws.Range("D2").Value = 41704
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 6
$ws.Range("C4").Value = 69
$ws.Range("D4").Value = 69
$ws.Range("C5").Value = 287
$ws.Range("D5").Value = 287
$ws.Range("C6").Value = 20129
$ws.Range("D6").Value = 20319
$ws.Range("C8").Value = 5292
$ws.Range("D8").Value = 5383
$ws.Range("C9").Value = 869
$ws.Range("D9").Value = 877
$ws.Range("C10").Value = 74
$ws.Range("D10").Value = 61
$ws.Range("C11").Value = 134
$ws.Range("D11").Value = 148
$ws.Range("C17").Value = 1744
$ws.Range("D17").Value = 1747
$ws.Range("C18").Value = 3988
$ws.Range("D18").Value = 3988
$ws.Range("C19").Value = 1441
$ws.Range("D19").Value = 1440
$ws.Range("C20").Value = 224460
$ws.Range("D20").Value = 224291
$ws.Range("C22").Value = 573
$ws.Range("D22").Value = 573
$ws.Range("C23").Value = 3959
$ws.Range("D23").Value = 3959
$ws.Range("C24").Value = 466
$ws.Range("D24").Value = 466
$ws.Range("C25").Value = 212893
$ws.Range("D25").Value = 212718
$ws.Range("C27").Value = 3353
$ws.Range("D27").Value = 3350
$ws.Range("C28").Value = 132
$ws.Range("D28").Value = 133
$ws.Range("C29").Value = 344
$ws.Range("D29").Value = 344
$ws.Range("C30").Value = 351
$ws.Range("D30").Value = 350
$ws.Range("C31").Value = 217
$ws.Range("D31").Value = 217
$ws.Range("C32").Value = 670
$ws.Range("D32").Value = 668
$ws.Range("C34").Value = 3189
$ws.Range("D34").Value = 3188
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 1
$ws.Range("C36").Value = 337
$ws.Range("D36").Value = 337
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("C38").Value = 71
$ws.Range("D38").Value = 71
$ws.Range("C39").Value = 621
$ws.Range("D39").Value = 620

# --- Fix D45 formula: was referencing D14, now should reference D18 ---
$ws.Range("D45").Formula = "=D18"

# --- Update active cell selection on the Table sheet ---
$ws.Range("E1").Select()

# --- Workbook view tab ratio change ---
$excel.ActiveWindow.TabRatio = 0.453

# --- Page setup: paper size 1 (Letter) -> 9 (A4) on every worksheet ---
foreach ($sheet in $wb.Worksheets) {
    $sheet.PageSetup.PaperSize = 9
}
